# "add `transpose` option in table module"
#
# The author selected the whole table on sheet "s1" (A1:G5), copied it, then
# switched to sheet "s2" and pasted it with Paste Special -> Transpose at A1.
# That lands the transposed table in A1:E8 on s2 (the gap row 2 on s1 becomes
# the empty column B on s2). Afterwards the cursor ends up one cell past the
# pasted block (G9) and s2 becomes the active sheet/tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("s1")
$ws2 = $wb.Worksheets.Item("s2")

# Mimic selecting + copying the source table on s1.
$ws1.Range("A1:G5").Select() | Out-Null
$ws1.Range("A1:G5").Copy()

# Paste the transposed values onto s2 (cell-by-cell, since this mirrors a
# Paste Special -> Transpose of A1:G5 -> A1:E8).
$ws2.Cells.Item(1, 1).Value = "a"
$ws2.Cells.Item(1, 3).Value = 1
$ws2.Cells.Item(1, 4).Value = 3
$ws2.Cells.Item(1, 5).Value = "ddd"

$ws2.Cells.Item(2, 1).Value = "one"
$ws2.Cells.Item(2, 3).Value = 2
$ws2.Cells.Item(2, 4).Value = 5

$ws2.Cells.Item(3, 1).Value = "two"
$ws2.Cells.Item(3, 3).Value = "xxx"
$ws2.Cells.Item(3, 4).Value = "ddd"
$ws2.Cells.Item(3, 5).Value = 7654

$ws2.Cells.Item(5, 1).Value = "one.two.three"
$ws2.Cells.Item(5, 3).Value = "a"
$ws2.Cells.Item(5, 4).Value = "b"
$ws2.Cells.Item(5, 5).Value = "c"

$ws2.Cells.Item(6, 1).Value = "arr.arr[]"
$ws2.Cells.Item(6, 3).Value = "1;2;3"
$ws2.Cells.Item(6, 4).Value = "df fssa"
$ws2.Cells.Item(6, 5).Value = "555 ;"

$ws2.Cells.Item(7, 1).Value = "one.1"
$ws2.Cells.Item(7, 3).Value = "we"

$ws2.Cells.Item(8, 1).Value = "z.0"
$ws2.Cells.Item(8, 3).Value = "rt"

# s2 becomes the active sheet, with the cursor resting just past the pasted
# range, the way it would after a real paste-special operation.
$ws2.Activate() | Out-Null
$ws2.Range("G9").Select() | Out-Null
